$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.406999999999999
$ws.Range("A3").Value = -21.972
$ws.Range("A14").Value = -21.659
$ws.Range("A16").Value = -22.051
$ws.Range("D18").Value = -8.225
$ws.Range("A21").Value = -19.947
$ws.Range("A23").Value = -20.302
$ws.Range("D24").Value = -7.574
$ws.Range("A25").Value = -21.759
$ws.Range("D25").Value = -7.709999999999999
$ws.Range("A26").Value = -21.435
$ws.Range("D27").Value = -8.279
$ws.Range("A29").Value = -21.232
$ws.Range("D30").Value = -7.305
$ws.Range("D31").Value = -7.930000000000001
$ws.Range("D39").Value = -7.502000000000001
$ws.Range("A40").Value = -20.182
$ws.Range("D48").Value = -7.100999999999999
$ws.Range("D51").Value = -8.259000000000002
$ws.Range("D52").Value = -8.199999999999999
$ws.Range("A53").Value = -21.924
$ws.Range("D55").Value = -8.258999999999999
$ws.Range("D56").Value = -8.459999999999999
$ws.Range("A57").Value = -22.576
$ws.Range("D57").Value = -8.270999999999999
$ws.Range("A59").Value = -22.571
$ws.Range("D60").Value = -8.440000000000001
$ws.Range("A65").Value = -21.482
$ws.Range("A69").Value = -21.587
$ws.Range("D73").Value = -8.004000000000001
$ws.Range("D74").Value = -7.964999999999999
$ws.Range("A79").Value = -20.856
$ws.Range("A83").Value = -22.035
$ws.Range("D89").Value = -6.822999999999999
$ws.Range("D90").Value = -7.629
$ws.Range("A91").Value = -21.509
$ws.Range("D92").Value = -6.621
$ws.Range("A93").Value = -21.395
$ws.Range("A100").Value = -21.95
